$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value. Every populated data row
# (2..206) had its "Förändrad" date bumped by exactly one day:
#   2023-10-03 (serial 45202) -> 2023-10-04 (serial 45203)
$firstRow = 2
$lastRow = 206

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45202) {
        $cell.Value = 45203
    }
}
